# Settings sheet: add a Name/Value/Description header table and one data
# row (QueueName -> YearlyReport), matching the new "Config" layout that
# replaced the previously-empty "Settings" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true

# Data row
$ws.Range("A2").Value = "QueueName"
$ws.Range("B2").Value = "YearlyReport"

# Column widths (closest values the host's pixel-grid ColumnWidth rounding
# can reach to the authored 28.109375 / 25.21875 / 45.44140625 char widths)
$ws.Columns.Item(1).ColumnWidth = 27.333333333333332
$ws.Columns.Item(2).ColumnWidth = 24.333333333333332
$ws.Columns.Item(3).ColumnWidth = 44.666666666666664

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it
$null = $ws.Range("C11").Select()
